# Update countries & provincias Spain
# Applies the COVID data refresh: updated stats for several countries
# (causing Maldivas to overtake Tanzania, and Guinea-Bisau to overtake
# Ruanda & Congo in the sorted list) plus a refreshed "last updated" time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 22:52"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1123764
$ws.Range("C4").Value = 28741
$ws.Range("E4").Value = 899386
$ws.Range("F4").Value = 16109
$ws.Range("G4").Value = 1601
$ws.Range("H4").Value = 65457

# --- Canada (row 15) ---
$ws.Range("B15").Value = 54810
$ws.Range("C15").Value = 1574
$ws.Range("D15").Value = 22515
$ws.Range("E15").Value = 28908

# --- Pakistan (row 27) ---
$ws.Range("B27").Value = 18092
$ws.Range("C27").Value = 1619
$ws.Range("E27").Value = 13324
$ws.Range("G27").Value = 56
$ws.Range("H27").Value = 417

# --- Tanzania / Maldivas swap around rows 113-114 ---
# Maldivas' updated numbers push it above Tanzania in the sorted table.
$ws.Range("A113").Value = "Maldivas"
$ws.Range("B113").Value = 491
$ws.Range("C113").Value = 23
$ws.Range("D113").Value = 17
$ws.Range("E113").Value = 473
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 1

$ws.Range("A114").Value = "Tanzania"
$ws.Range("B114").Value = 480
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 167
$ws.Range("E114").Value = 297
$ws.Range("F114").Value = 7
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 16

# --- Guinea-Bisau / Ruanda / Congo reorder around rows 132-134 ---
# Guinea-Bisau's updated numbers push it above Ruanda and Congo.
$ws.Range("A132").Value = "Guinea-Bisau"
$ws.Range("B132").Value = 257
$ws.Range("C132").Value = 52
$ws.Range("D132").Value = 19
$ws.Range("E132").Value = 237
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 1

$ws.Range("A133").Value = "Ruanda"
$ws.Range("B133").Value = 243
$ws.Range("C133").Value = 0
$ws.Range("D133").Value = 104
$ws.Range("E133").Value = 139
$ws.Range("F133").Value = 0
$ws.Range("G133").Value = 0
$ws.Range("H133").Value = 0

$ws.Range("A134").Value = "Congo"
$ws.Range("B134").Value = 220
$ws.Range("C134").Value = 0
$ws.Range("D134").Value = 19
$ws.Range("E134").Value = 192
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 9

# --- Trinidad yTobago (row 149) ---
$ws.Range("D149").Value = 81
$ws.Range("E149").Value = 27
